# Regenerate save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" column (column G) values for rows 2-28, replacing the old Strike# values
$kValues = @(3, 3, 0, 6, 1, 2, 8, 7, 6, 6, 4, 5, 1, 0, 6, 3, 5, 4, 6, 7, 3, 3, 3, 2, 4, 1, 1)

$startRow = 2
for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
